$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 374
$ws1.Range("F6").Value = 2027
$ws1.Range("F7").Value = 107

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 108

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 374
$ws4.Range("F6").Value = 108
$ws4.Range("F10").Value = 2027
$ws4.Range("F11").Value = 107
